$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.165.91'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +7.76%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.586.02'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +7.65%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9904'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.93%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '298.32'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.48%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3616'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3334'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +8.31%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '41.27'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.112'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.34%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06929'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.08%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.13%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.37'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.13%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.812'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.30%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.505'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.39%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9911'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.00%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001060'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.582.76'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.40%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06592'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +10.95%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '75.87'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +10.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.77'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.28%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.901'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.56'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.184.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.88%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.355'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.33%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.492'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +16.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '148.29'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.04%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +11.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.755.86'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.45%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.08%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.922'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.81%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.859'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +18.16%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9157'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +13.61%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08129'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.628'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.53%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.66'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +11.79%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.102'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.234'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06013'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.23%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.266'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.54%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02178'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.76%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1972'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9906'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.89%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5759'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.84%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.763'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.86%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.77'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.16%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '124.60'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5548'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.43%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.931'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.98%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06702'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.62%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.01'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.96%  '

